$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.702.95"
$ws.Range("E2").Value = "  +2.67%  "

$ws.Range("D3").Value = "4.034.96"
$ws.Range("E3").Value = "  +2.20%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("E7").Value = "  +1.13%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.743"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.79%  "

$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.26%  "

$ws.Range("D14").Value = "4.671.92"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D15").Value = "4.040.61"
$ws.Range("E15").Value = "  +2.39%  "

$ws.Range("E16").Value = "  +7.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("D20").Value = "71.688.44"
$ws.Range("E20").Value = "  +2.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "700.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("E31").Value = "  +3.20%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "68.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.78%  "

$ws.Range("D35").Value = "0.0₃0915"
$ws.Range("E35").Value = "  +4.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.448"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.50%  "

$ws.Range("E38").Value = "  +3.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.06%  "

$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.20%  "

$ws.Range("E46").Value = "  +2.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000279"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("D51").Value = "0.0₆0346"
$ws.Range("E51").Value = "  -4.79%  "
